$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the attribute values: Color row gains "Yellow", Legs row becomes material text
$ws.Range("H2").Value = "White,Yellow"
$ws.Range("H3").Value = "Aluminium,Steel"

# Move the active selection to A2 (was I10)
$ws.Range("A2").Select()
